$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# ------------------------------------------------------------------
# Move the existing "En negro / En azul" note from row 31 down to
# row 37, making room for the new rows (27-29) below it.
# ------------------------------------------------------------------
$ws.Range("A31").Cut($ws.Range("A37"))

# ------------------------------------------------------------------
# Fill in the new content (adds 3 new shared strings):
#   B27 -> IndiceAridad
#   B28 -> IndiceContinentalidad
#   A29 -> "the rest follow the schema on the otalex config file...."
# ------------------------------------------------------------------
$ws.Range("B27").Value2 = "IndiceAridad"
$ws.Range("B28").Value2 = "IndiceContinentalidad"
$ws.Range("A29").Value2 = "the rest follow the schema on the otalex config file. This will have to be completed from there."

# Give the new cells the same style used throughout this block
# (the "s=3" / blue comment style already used by row 26, e.g. A26).
$ws.Range("A26").Copy()
$ws.Range("B27").PasteSpecial(-4122)
$ws.Range("B28").PasteSpecial(-4122)
$ws.Range("A29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Update the view / selection to match the edited area.
# ------------------------------------------------------------------
[void]$ws.Range("A29").Select()
